$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functional Requirements Sheet")
$ws.Range("G3:G12").Validation.Delete()
$ws.Range("G17").Validation.Delete()
$ws.Range("G3:G12").Validation.Add(3, 1, 1, "Mobile App,Online Banking Website,Both")
$ws.Range("G17").Validation.Add(3, 1, 1, "Mobile App,Online Banking Website,Both")
Write-Host "ok"
